# Update the "CCmd" sheet: recompute the LJ-optimized (col G) and
# X6-optimized (col I) series, then leave the selection on the newly
# plotted range (G2:G16), matching what Excel leaves selected after the
# user charts that range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCmd")

# New values for column G (LJ-optimized) and column I (X6-optimized), rows 2-16.
$values = @{
    2  = @{ G = -1.1224000000000001; I = -1.1273 }
    3  = @{ G = -1.1551;              I = -1.1593 }
    4  = @{ G = -1.1852;              I = -1.1886000000000001 }
    5  = @{ G = -1.2121999999999999; I = -1.2143999999999999 }
    6  = @{ G = -1.2352000000000001; I = -1.2361 }
    7  = @{ G = -1.2532000000000001; I = -1.2525999999999999 }
    8  = @{ G = -1.2652000000000001; I = -1.2632000000000001 }
    9  = @{ G = -1.27;                I = -1.2665 }
    10 = @{ G = -1.266;               I = -1.2613000000000001 }
    11 = @{ G = -1.2516;              I = -1.2462 }
    12 = @{ G = -1.2246999999999999; I = -1.2194 }
    13 = @{ G = -1.1832;              I = -1.1791 }
    14 = @{ G = -1.1243000000000001; I = -1.1232 }
    15 = @{ G = -1.0448999999999999; I = -1.0491999999999999 }
    16 = @{ G = -0.94140000000000001; I = -0.95450000000000002 }
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row].G
    $ws.Cells.Item($row, 9).Value = $values[$row].I
}

# Make CCmd the active sheet and select the range that was just plotted
# (G2:G16), matching the saved selection state in the workbook.
$ws.Activate()
$ws.Range("G2:G16").Select()
